$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows above the existing 2014 data row so it becomes row 4,
#     leaving room for the 2012 and 2013 rows to be inserted at rows 2 and 3.
$ws.Rows("2:3").Insert()
# Row insert picks up the header row's formatting by default; the new data
# rows should carry the plain/default style (like the original data row did).
$ws.Rows("2:3").Style = "Normal"

# --- Insert new blank columns so each of the 4 split metrics gets a "_y" column
#     right after its "_x" column. Insert from rightmost to leftmost so the
#     column letters used below still refer to the intended original columns.
#     After all 4 inserts:
#       A=code B=year C=crop D=temp E=agri_x F=agri_y(NEW) G=precip
#       H=fert_x I=fert_y(NEW) J=pop_x K=pop_y(NEW) L=emp_x M=emp_y(NEW)
$ws.Columns("J").Insert()
$ws.Columns("I").Insert()
$ws.Columns("H").Insert()
$ws.Columns("F").Insert()

# --- Header row (row 1): rename split headers and fill the new "_y" columns ---
$ws.Range("E1").Value = "4. Agriculture land area (% of land area)_x"
$ws.Range("F1").Value = "4. Agriculture land area (% of land area)_y"
$ws.Range("H1").Value = "7. Fertilizer consumption (kilograms per hectare of arable land)_x"
$ws.Range("I1").Value = "7. Fertilizer consumption (kilograms per hectare of arable land)_y"
$ws.Range("J1").Value = "13. Population_x"
$ws.Range("K1").Value = "13. Population_y"
$ws.Range("L1").Value = "17. Employment in agriculture (% of total employment) (modeled ILO estimate)_x"
$ws.Range("M1").Value = "17. Employment in agriculture (% of total employment) (modeled ILO estimate)_y"

# The new header cells (F1,I1,K1,M1) need the same bold/bordered/centered style
# as the rest of row 1 -- set the formatting directly (Style object assignment
# doesn't propagate through this host, so set the underlying attributes).
foreach ($addr in @("F1", "I1", "K1", "M1")) {
    $c = $ws.Range($addr)
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
    $c.Borders.LineStyle = 1
}

# --- Duplicate the "_x" values into the new "_y" columns for the pre-existing 2014 row (now row 4) ---
$ws.Range("F4").Value = $ws.Range("E4").Value2
$ws.Range("I4").Value = $ws.Range("H4").Value2
$ws.Range("K4").Value = $ws.Range("J4").Value2
$ws.Range("M4").Value = $ws.Range("L4").Value2

# --- New data row for 2012 (row 2) ---
$ws.Range("A2").Value = "SDN"
$ws.Range("B2").Value = 2012
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "69.56"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = 27.55
$ws.Range("E2").Value = 59.6712848
$ws.Range("F2").Value = 59.6712848
$ws.Range("G2").Value = 253.3
$ws.Range("H2").Value = 3.327887329
$ws.Range("I2").Value = 3.327887329
$ws.Range("J2").Value = 35159792
$ws.Range("K2").Value = 35159792
$ws.Range("L2").Value = 44.549260091789
$ws.Range("M2").Value = 44.549260091789

# --- New data row for 2013 (row 3) ---
$ws.Range("A3").Value = "SDN"
$ws.Range("B3").Value = 2013
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "91.93"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = 27.54
$ws.Range("E3").Value = 59.66453961
$ws.Range("F3").Value = 59.66453961
$ws.Range("G3").Value = 241.01
$ws.Range("H3").Value = 2.612458392
$ws.Range("I3").Value = 2.612458392
$ws.Range("J3").Value = 35990704
$ws.Range("K3").Value = 35990704
$ws.Range("L3").Value = 44.037527119395
$ws.Range("M3").Value = 44.037527119395

$wb.Save()
